$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.219.51"
$ws.Range("E2").Value = "  +0.18%  "
$ws.Range("D3").Value = "1.860.59"
$ws.Range("E3").Value = "  +0.74%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "'0.7139"
$ws.Range("E5").Value = "  +1.62%  "
$ws.Range("D6").Value = "'238.03"
$ws.Range("E6").Value = "  -0.21%  "
$ws.Range("D8").Value = "'0.08185"
$ws.Range("E8").Value = "  +10.43%  "
$ws.Range("D9").Value = "'0.3049"
$ws.Range("E9").Value = "  +0.02%  "
$ws.Range("D10").Value = "'23.22"
$ws.Range("E10").Value = "  -0.69%  "
$ws.Range("D11").Value = "'0.08175"
$ws.Range("E11").Value = "  +0.52%  "
$ws.Range("D12").Value = "1.910.81"
$ws.Range("E12").Value = "  +1.30%  "
$ws.Range("D13").Value = "'5.175"
$ws.Range("E13").Value = "  -0.73%  "
$ws.Range("D15").Value = "'89.73"
$ws.Range("E15").Value = "  +1.03%  "
$ws.Range("D16").Value = "29.230.16"
$ws.Range("E16").Value = "  -0.17%  "
$ws.Range("D17").Value = "'0.000007940"
$ws.Range("E17").Value = "  +3.84%  "
$ws.Range("D18").Value = "'5.794"
$ws.Range("E18").Value = "  +0.36%  "
$ws.Range("D19").Value = "'13.37"
$ws.Range("E19").Value = "  +2.34%  "
$ws.Range("D20").Value = "'237.20"
$ws.Range("E20").Value = "  -0.48%  "
$ws.Range("E21").Value = "  +0.00%  "
$ws.Range("D22").Value = "2.109.86"
$ws.Range("E22").Value = "  -0.29%  "
$ws.Range("D23").Value = "'1.001"
$ws.Range("E23").Value = "  -0.04%  "
$ws.Range("D24").Value = "'7.433"
$ws.Range("E24").Value = "  -2.11%  "
$ws.Range("D25").Value = "'162.69"
$ws.Range("E25").Value = "  +1.22%  "
$ws.Range("B26").Value = "Stellar"
$ws.Range("C26").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D26").Value = "'0.1467"
$ws.Range("E26").Value = "  +1.12%  "
$ws.Range("B27").Value = "Cosmos"
$ws.Range("C27").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D27").Value = "'8.963"
$ws.Range("E27").Value = "  -0.41%  "
$ws.Range("D28").Value = "'18.10"
$ws.Range("E28").Value = "  +0.06%  "
$ws.Range("D29").Value = "'1.959"
$ws.Range("E29").Value = "  -0.50%  "
$ws.Range("E30").Value = "  +1.87%  "
$ws.Range("D31").Value = "'1.484"
$ws.Range("E31").Value = "  -0.37%  "
$ws.Range("E32").Value = "  -2.62%  "
$ws.Range("E33").Value = "  +0.81%  "
$ws.Range("D34").Value = "'0.05226"
$ws.Range("E34").Value = "  +0.85%  "
$ws.Range("E35").Value = "  -1.19%  "
$ws.Range("D36").Value = "'0.7086"
$ws.Range("E36").Value = "  +0.77%  "
$ws.Range("E37").Value = "  -2.90%  "
$ws.Range("D38").Value = "'2.674"
$ws.Range("E38").Value = "  +0.24%  "
$ws.Range("D39").Value = "'0.01861"
$ws.Range("E39").Value = "  -0.40%  "
$ws.Range("D40").Value = "'2.728"
$ws.Range("E40").Value = "  +2.03%  "
$ws.Range("D41").Value = "'0.9236"
$ws.Range("E41").Value = "  -2.44%  "
$ws.Range("D42").Value = "1.141.56"
$ws.Range("E42").Value = "  +7.05%  "
$ws.Range("D43").Value = "'0.4287"
$ws.Range("E43").Value = "  -0.29%  "
$ws.Range("D44").Value = "'5.898"
$ws.Range("E44").Value = "  -1.83%  "
$ws.Range("D45").Value = "'70.45"
$ws.Range("E45").Value = "  -0.26%  "
$ws.Range("D46").Value = "'0.9998"
$ws.Range("E46").Value = "  -0.05%  "
$ws.Range("E47").Value = "  -0.11%  "
$ws.Range("D48").Value = "'1.776"
$ws.Range("E48").Value = "  +1.99%  "
$ws.Range("D49").Value = "2.007.63"
$ws.Range("E49").Value = "  +0.06%  "
$ws.Range("D50").Value = "'9.211"
$ws.Range("E50").Value = "  +1.19%  "
$ws.Range("D51").Value = "'6.959"
$ws.Range("E51").Value = "  -1.00%  "
